$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation is inserted as the new row 44, pushing the
# existing rows 44:75 down to 45:76 (dimension grows from A1:R75 to A1:R76).
$ws.Rows("44:44").Insert()

$ws.Cells.Item(44, 1).Value  = 7
$ws.Cells.Item(44, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(44, 3).Value  = "Ñuble"
$ws.Cells.Item(44, 4).Value  = "2022-03-31"
$ws.Cells.Item(44, 5).Value  = 16
$ws.Cells.Item(44, 6).Value  = 100112021
$ws.Cells.Item(44, 7).Value  = "Ají"
$ws.Cells.Item(44, 8).Value  = "Americana (o)"
$ws.Cells.Item(44, 9).Value  = "Primera"
$ws.Cells.Item(44, 10).Value = 60
$ws.Cells.Item(44, 11).Value = 8500
$ws.Cells.Item(44, 12).Value = 9000
$ws.Cells.Item(44, 13).Value = 8750
$ws.Cells.Item(44, 14).Value = "`$/caja 15 kilos"
$ws.Cells.Item(44, 15).Value = "Región del Maule"
$ws.Cells.Item(44, 16).Value = 583
$ws.Cells.Item(44, 17).Value = 15
$ws.Cells.Item(44, 18).Value = "Hortaliza"
